$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "level 1" dialog / instruction rows (rows 20-29, columns A=Key, B=Value).
# Values are written in the same order the original author entered them so the
# shared-string table is rebuilt in the matching sequence.
$ws.Range("A20").Value = "level_1_intro_1"
$ws.Range("B20").Value = "A chasm impedes Pengu's journey."
$ws.Range("A21").Value = "level_1_intro_2"
$ws.Range("B21").Value = "Help Pengu out by determining the distance of the gap!"
$ws.Range("A22").Value = "level_1_info_1"
$ws.Range("A26").Value = "level_1_info_3"
$ws.Range("A27").Value = "level_1_info_4"
$ws.Range("A28").Value = "level_1_info_5"
$ws.Range("A23").Value = "level_1_info_2_a"
$ws.Range("A24").Value = "level_1_info_2_b"
$ws.Range("A29").Value = "success"
$ws.Range("B29").Value = "SUCCESS!"
$ws.Range("B22").Value = "In order to add these two fractions, you must make the denominators equal."
$ws.Range("B23").Value = "One way to do this is by multiplying the denominators together to make them equal."
$ws.Range("B24").Value = "Then multiply the numerators by the same amount from the denominator."
$ws.Range("A25").Value = "level_1_info_2_c"
$ws.Range("B25").Value = "After that, you can add both fractions properly."
$ws.Range("B26").Value = "Now it’s your turn! Use the multiplier to make both denominators equal."
$ws.Range("B28").Value = "Go ahead and type in the correct answer by pressing on either slot!"
$ws.Range("B27").Value = "Notice how both fractions now have the same unit sizes? They can now be added properly."

# B28 gets a vertically-centered style (new cellXfs entry).
$ws.Range("B28").VerticalAlignment = -4108

# Leave the selection on B28, matching the saved view state.
$ws.Range("B28").Select()
